$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the diff (D = Price, E = Volume(1h)).
# Cells whose new text looks like a plain number get NumberFormat "@" forced
# first so Excel keeps them as text (matching the original inlineStr type)
# instead of silently re-parsing them into a numeric value.

$ws.Range('D2').Value = '54.844.41'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '2.271.02'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '504.49'
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '127.95'
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  +1.03%  '
$ws.Range('D9').Value = '2.281.78'
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0981'
$ws.Range('E10').Value = '  +3.32%  '
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.08'
$ws.Range('E12').Value = '  +7.54%  '
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '23.53'
$ws.Range('E14').Value = '  +4.37%  '
$ws.Range('D15').Value = '2.674.35'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '55.007.92'
$ws.Range('E16').Value = '  +1.97%  '
$ws.Range('E17').Value = '  +1.11%  '
$ws.Range('D18').Value = '2.286.15'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.35'
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('E20').Value = '  +1.49%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '313.59'
$ws.Range('E21').Value = '  +4.31%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.56'
$ws.Range('E22').Value = '  +4.32%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '59.89'
$ws.Range('E24').Value = '  -1.43%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.996'
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('E26').Value = '  +4.20%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.52'
$ws.Range('E27').Value = '  +3.68%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '171.17'
$ws.Range('E29').Value = '  +4.12%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.64'
$ws.Range('E30').Value = '  +2.42%  '
$ws.Range('E31').Value = '  +2.38%  '
$ws.Range('E32').Value = '  +7.17%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '17.96'
$ws.Range('E34').Value = '  +1.58%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.995'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.23'
$ws.Range('E36').Value = '  +3.26%  '
$ws.Range('E37').Value = '  -3.11%  '
$ws.Range('E38').Value = '  +5.36%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.67'
$ws.Range('E39').Value = '  +2.45%  '
$ws.Range('E40').Value = '  +4.48%  '
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '135.74'
$ws.Range('E42').Value = '  +8.97%  '
$ws.Range('E43').Value = '  +3.94%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.87'
$ws.Range('E44').Value = '  +2.11%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '258.37'
$ws.Range('E45').Value = '  +9.26%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0504'
$ws.Range('E46').Value = '  +3.09%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0919'
$ws.Range('E47').Value = '  +3.70%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.547'
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('E49').Value = '  +4.14%  '
$ws.Range('E50').Value = '  +1.16%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '16.44'
$ws.Range('E51').Value = '  +2.31%  '
